$wb = $excel.ActiveWorkbook

# The "Charts" tab (5th sheet) currently has no content and is not the
# active tab. Add the "coming soon" placeholder message to it and make
# it the active/selected sheet, matching the author's intent of adding a
# coming-soon message to the charts tab of the generated report.
$ws = $wb.Worksheets.Item("Charts")
$ws.Range("A1").Value = "Automatically generated chart(s) coming soon to this tab."
$ws.Activate()
